$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# --- Update the per-row "time_taken" query timestamps in column F (rows 2-19) ---
$timestamps = @(
    "2021-10-05 14:35:57.175360",
    "2021-10-05 14:35:57.175368",
    "2021-10-05 14:35:57.175371",
    "2021-10-05 14:35:57.175374",
    "2021-10-05 14:35:57.175376",
    "2021-10-05 14:35:57.175379",
    "2021-10-05 14:35:57.175382",
    "2021-10-05 14:35:57.175384",
    "2021-10-05 14:35:57.175387",
    "2021-10-05 14:35:57.175389",
    "2021-10-05 14:35:57.175392",
    "2021-10-05 14:35:57.175394",
    "2021-10-05 14:35:57.175397",
    "2021-10-05 14:35:57.175399",
    "2021-10-05 14:35:57.175402",
    "2021-10-05 14:35:57.175405",
    "2021-10-05 14:35:57.175408",
    "2021-10-05 14:35:57.175410"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}

# --- Add the new "metadata" worksheet as a second tab, right after "data" ---
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

# Header row (B1:G1): reuse the bold/centered/bordered header style already
# used on the "data" sheet's header row.
$ws.Range("B1:F1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)
$meta.Cells.Item(1, 2).Value = "data_name"
$meta.Cells.Item(1, 3).Value = "data_id"
$meta.Cells.Item(1, 4).Value = "data_version"
$meta.Cells.Item(1, 5).Value = "data_version_created"
$meta.Cells.Item(1, 6).Value = "panel_query_time"
$meta.Cells.Item(1, 7).Value = "panel_get_request"

# A2: numeric row index, styled like the "data" sheet's A column entries.
$ws.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$meta.Cells.Item(2, 1).Value = 0

$meta.Cells.Item(2, 2).Value = "Vitreoretinopathy"
$meta.Cells.Item(2, 3).Value = 3113

# data_version ("1.2") must be stored as literal text, not a number. Write it
# as a text formula and immediately flatten it to a static value in place so
# it lands as a plain string without picking up a new number-format style.
$d2 = $meta.Cells.Item(2, 4)
$d2.Formula = '="1.2"'
$d2.Copy()
$d2.PasteSpecial(-4163)

$meta.Cells.Item(2, 5).Value = "2021-06-07T06:51:41.373228Z"
$meta.Cells.Item(2, 6).Value = "2021-10-05 14:35:57.171854"
$meta.Cells.Item(2, 7).Value = "https://panelapp.agha.umccr.org/api/v1/panels/3113/?format=json"

# Keep "data" as the active tab (matches the original workbook's activeTab=0,
# which the commit's diff leaves untouched).
$ws.Activate()
